# Update BOC USD rates (auto)
# Adds a new published-rate row to "All Published Values" (row 7),
# grows the AutoFilter/_FilterDatabase range to include it, and bumps
# the "publishes" count for 2026-01-02 on the "Daily Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Published Values" -----------------------------------
$ws1 = $wb.Worksheets.Item("All Published Values")

# Format the new row as text first so the values round-trip as literal
# strings (matching the existing rows) instead of being auto-converted
# to dates/numbers by Excel.
$ws1.Range("A7:J7").NumberFormat = "@"

$ws1.Range("A7").Value = "2026-01-02"
$ws1.Range("B7").Value = "2026-01-02 19:23:08"
$ws1.Range("C7").Value = "697.85"
$ws1.Range("D7").Value = "697.85"
$ws1.Range("E7").Value = "700.79"
$ws1.Range("F7").Value = "700.79"
$ws1.Range("G7").Value = "702.88"
$ws1.Range("H7").Value = "2026/01/02 19:23:08"
$ws1.Range("I7").Value = "2026-01-02 11:28:24"
$ws1.Range("J7").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Remove the temporary text formatting so the new row has no explicit
# style, same as the other data rows.
$ws1.Range("A7:J7").ClearFormats()

# Expand the AutoFilter range to cover the new row (toggle off/on so the
# new range actually takes effect instead of just removing the filter).
$ws1.AutoFilterMode = $false
$ws1.Range("A1:J7").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$7"
    }
}

# --- Sheet 2: "Daily Summary" -------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B4").Value = 6
